$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.771.91'

$ws.Range("D3").Value = '1.796.64'

$ws.Range("D4").Value = '0.9995'

$ws.Range("D5").Value = '309.23'
$ws.Range("E5").Value = '  -0.55%  '

$ws.Range("D6").Value = '0.9998'
$ws.Range("E6").Value = '  -0.09%  '

$ws.Range("D7").Value = '0.4407'
$ws.Range("E7").Value = '  +4.41%  '

$ws.Range("D8").Value = '0.3669'
$ws.Range("E8").Value = '  -0.25%  '

$ws.Range("D9").Value = '0.07315'
$ws.Range("E9").Value = '  +1.48%  '

$ws.Range("D10").Value = '0.8533'
$ws.Range("E10").Value = '  +1.62%  '

$ws.Range("B11").Value = 'WrappedEther'
$ws.Range("C11").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D11").Value = '1.981.47'
$ws.Range("E11").Value = '  +8.85%  '

$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D12").Value = '20.58'
$ws.Range("E12").Value = '  -1.08%  '

$ws.Range("D13").Value = '6.601'
$ws.Range("E13").Value = '  -0.78%  '

$ws.Range("E14").Value = '  +2.08%  '

$ws.Range("D15").Value = '0.07042'
$ws.Range("E15").Value = '  -0.14%  '

$ws.Range("D16").Value = '5.257'
$ws.Range("E16").Value = '  -0.34%  '

$ws.Range("D17").Value = '0.9999'
$ws.Range("E17").Value = '  -0.14%  '

$ws.Range("D18").Value = '0.000008643'
$ws.Range("E18").Value = '  -1.51%  '

$ws.Range("D19").Value = '0.9995'
$ws.Range("E19").Value = '  -0.06%  '

$ws.Range("E20").Value = '  -1.15%  '

$ws.Range("D21").Value = '26.814.82'
$ws.Range("E21").Value = '  -1.12%  '

$ws.Range("D22").Value = '5.138'
$ws.Range("E22").Value = '  +0.26%  '

$ws.Range("E23").Value = '  -0.52%  '

$ws.Range("D24").Value = '1.973'
$ws.Range("E24").Value = '  -0.05%  '

$ws.Range("D25").Value = '151.49'
$ws.Range("E25").Value = '  -0.24%  '

$ws.Range("D26").Value = '2.195'
$ws.Range("E26").Value = '  -1.66%  '

$ws.Range("D27").Value = '18.31'
$ws.Range("E27").Value = '  +0.35%  '

$ws.Range("D28").Value = '5.173'
$ws.Range("E28").Value = '  -1.71%  '

$ws.Range("D29").Value = '116.81'
$ws.Range("E29").Value = '  +0.62%  '

$ws.Range("D30").Value = '0.08774'
$ws.Range("E30").Value = '  +0.30%  '

$ws.Range("D31").Value = '0.7362'
$ws.Range("E31").Value = '  -0.34%  '

$ws.Range("D32").Value = '1.151'
$ws.Range("E32").Value = '  -2.05%  '

$ws.Range("D33").Value = '2.902'
$ws.Range("E33").Value = '  -0.42%  '

$ws.Range("D34").Value = '4.417'
$ws.Range("E34").Value = '  +0.02%  '

$ws.Range("D35").Value = '0.9991'
$ws.Range("E35").Value = '  -0.10%  '

$ws.Range("E36").Value = '  -0.68%  '

$ws.Range("D37").Value = '0.01954'
$ws.Range("E37").Value = '  +0.11%  '

$ws.Range("D38").Value = '0.05179'
$ws.Range("E38").Value = '  -1.28%  '

$ws.Range("D39").Value = '0.5208'
$ws.Range("E39").Value = '  +3.48%  '

$ws.Range("D40").Value = '7.026'
$ws.Range("E40").Value = '  -4.30%  '

$ws.Range("E41").Value = '  -2.70%  '

$ws.Range("D42").Value = '0.1674'
$ws.Range("E42").Value = '  -0.86%  '

$ws.Range("D43").Value = '0.4987'
$ws.Range("E43").Value = '  +5.92%  '

$ws.Range("D44").Value = '8.403'
$ws.Range("E44").Value = '  -1.97%  '

$ws.Range("E45").Value = '  +4.46%  '

$ws.Range("E46").Value = '  -1.53%  '

$ws.Range("D47").Value = '104.92'
$ws.Range("E47").Value = '  -1.22%  '

$ws.Range("D48").Value = '0.9989'
$ws.Range("E48").Value = '  -0.10%  '

$ws.Range("E49").Value = '  +0.55%  '

$ws.Range("E50").Value = '  -0.66%  '

$ws.Range("D51").Value = '0.9146'
$ws.Range("E51").Value = '  +1.58%  '
